$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 178
$ws.Range("G3").Value = "不可售"
$ws.Range("F4").Value = 136
$ws.Range("F5").Value = 1291
$ws.Range("F6").Value = 18014
$ws.Range("F7").Value = 356
$ws.Range("F9").Value = 1067
$ws.Range("F10").Value = 6797
$ws.Range("F11").Value = 683
$ws.Range("F12").Value = 155
$ws.Range("F13").Value = 11
$ws.Range("F14").Value = 109
$ws.Range("F19").Value = 210
$ws.Range("F21").Value = 653
$ws.Range("F25").Value = 271
$ws.Range("F26").Value = 980
$ws.Range("F27").Value = 113
$ws.Range("F28").Value = 5161
$ws.Range("F29").Value = 532
$ws.Range("F30").Value = 23
$ws.Range("F31").Value = 13
$ws.Range("F33").Value = 12020
$ws.Range("F34").Value = 1275
$ws.Range("F35").Value = 40
$ws.Range("F36").Value = 203
$ws.Range("F38").Value = 3910

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 2

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 178
$ws.Range("G3").Value = "不可售"
$ws.Range("F4").Value = 136
$ws.Range("F5").Value = 1291
$ws.Range("F6").Value = 18014
$ws.Range("F7").Value = 356
$ws.Range("F9").Value = 1067
$ws.Range("F10").Value = 6797
$ws.Range("F11").Value = 683
$ws.Range("F12").Value = 155
$ws.Range("F13").Value = 11
$ws.Range("F14").Value = 109
$ws.Range("F19").Value = 210
$ws.Range("F21").Value = 653
$ws.Range("F25").Value = 271
$ws.Range("F26").Value = 980
$ws.Range("F27").Value = 113
$ws.Range("F28").Value = 5161
$ws.Range("F29").Value = 532
$ws.Range("F30").Value = 2
$ws.Range("F32").Value = 23
$ws.Range("F33").Value = 13
$ws.Range("F35").Value = 12020
$ws.Range("F36").Value = 1275
$ws.Range("F37").Value = 40
$ws.Range("F38").Value = 203
$ws.Range("F40").Value = 3910
